$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.203398585319519
$ws.Range("B1").Value = 1.836362838745117
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.848986268043518
$ws.Range("E1").Value = 1.205420017242432
